{"js": "// 03 classification validation & restructure future labs\n// Update generated timestamp and rotate sample customer names/SSNs.\n\nconst replacements = [\n  [\"Generated: 2025-12-30 09:36\", \"Generated: 2025-12-31 08:15\"],\n  [\"Amanda Rodriguez\", \"Lauren Smith\"],\n  [\"021-08-2161\", \"605-50-4048\"],\n  [\"Michelle Moore\", \"Robert Rodriguez\"],\n  [\"026-35-7420\", \"339-93-9227\"],\n  [\"Daniel Williams\", \"Lauren Smith\"],\n  [\"849-74-2296\", \"583-06-4554\"],\n  [\"William Wilson\", \"Amanda Wilson\"],\n  [\"528-89-0681\", \"475-47-0764\"],\n  [\"William Taylor\", \"Jessica Moore\"],\n  [\"279-85-6302\", \"854-45-5056\"],\n  [\"Jennifer Martin\", \"Lauren Davis\"],\n  [\"336-31-0519\", \"714-18-1324\"],\n  [\"Michael Moore\", \"Emily Davis\"],\n  [\"455-89-1533\", \"169-77-8145\"],\n  [\"Lisa Wilson\", \"Jennifer Gonzalez\"],\n  [\"746-53-8112\", \"662-48-7773\"],\n];\n\nconst body = context.document.body;\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# 03 classification validation & restructure future labs\n# Update generated timestamp and rotate sample customer names/SSNs.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-Text \"Generated: 2025-12-30 09:36\" \"Generated: 2025-12-31 08:15\"\n\nReplace-Text \"Amanda Rodriguez\" \"Lauren Smith\"\nReplace-Text \"021-08-2161\" \"605-50-4048\"\n\nReplace-Text \"Michelle Moore\" \"Robert Rodriguez\"\nReplace-Text \"026-35-7420\" \"339-93-9227\"\n\nReplace-Text \"Daniel Williams\" \"Lauren Smith\"\nReplace-Text \"849-74-2296\" \"583-06-4554\"\n\nReplace-Text \"William Wilson\" \"Amanda Wilson\"\nReplace-Text \"528-89-0681\" \"475-47-0764\"\n\nReplace-Text \"William Taylor\" \"Jessica Moore\"\nReplace-Text \"279-85-6302\" \"854-45-5056\"\n\nReplace-Text \"Jennifer Martin\" \"Lauren Davis\"\nReplace-Text \"336-31-0519\" \"714-18-1324\"\n\nReplace-Text \"Michael Moore\" \"Emily Davis\"\nReplace-Text \"455-89-1533\" \"169-77-8145\"\n\nReplace-Text \"Lisa Wilson\" \"Jennifer Gonzalez\"\nReplace-Text \"746-53-8112\" \"662-48-7773\"\n"}
